# Round formula results to 2 decimal places
#
# Wraps the existing multiplication / MIN formulas in several blocks of the
# "Presets" sheet with ROUND(...,2) so displayed/stored results are rounded
# to 2 decimal places.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("C", "D", "E", "F", "G")

# Block 1: rows 5-11, base multiplier row is row 3 (e.g. C5 = B5*C$3)
foreach ($row in 5..11) {
    foreach ($col in $cols) {
        $ws.Range("$col$row").Formula = "=ROUND(B$row*$col`$3,2)"
    }
}

# Block 2: rows 18-24, base multiplier row is row 16 (e.g. C18 = B18*C$16)
foreach ($row in 18..24) {
    foreach ($col in $cols) {
        $ws.Range("$col$row").Formula = "=ROUND(B$row*$col`$16,2)"
    }
}

# Block 3: rows 31-37, base multiplier row is row 29 (e.g. C31 = B31*C$29)
foreach ($row in 31..37) {
    foreach ($col in $cols) {
        $ws.Range("$col$row").Formula = "=ROUND(B$row*$col`$29,2)"
    }
}

# Block 4: rows 44-50, columns D-G only (C holds a literal 1, untouched),
# base multiplier row is row 42 (e.g. D44 = MIN(B44*D$42,1))
foreach ($row in 44..50) {
    foreach ($col in @("D", "E", "F", "G")) {
        $ws.Range("$col$row").Formula = "=ROUND(MIN(B$row*$col`$42,1),2)"
    }
}
